$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.8
$ws.Range("I2").Value = 2.63
$ws.Range("X2").Value = 13
$ws.Range("AJ2").Value = 26
$ws.Range("AN2").Value = 4.75
$ws.Range("BA2").Value = 81

# Row 5
$ws.Range("S5").Value = 1.62

# Row 6
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("S6").Value = 1.62

# Row 7
$ws.Range("S7").Value = 1.67

# Row 8
$ws.Range("S8").Value = 1.36

# Row 10
$ws.Range("G10").Value = 3.25
$ws.Range("I10").Value = 2.25
$ws.Range("W10").Value = 8.5
$ws.Range("AK10").Value = 21
$ws.Range("AW10").Value = 4.33

# Row 11
$ws.Range("G11").Value = 2.1
$ws.Range("I11").Value = 3.8
$ws.Range("J11").Value = 2.88
$ws.Range("L11").Value = 4.33
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("X11").Value = 9.5
$ws.Range("Y11").Value = 9.5
$ws.Range("Z11").Value = 19
$ws.Range("AA11").Value = 19
$ws.Range("AH11").Value = 19
$ws.Range("AI11").Value = 13
$ws.Range("AO11").Value = 12
$ws.Range("AX11").Value = 21
$ws.Range("AZ11").Value = 67
$ws.Range("BB11").Value = 251
